# "changes in config file done from my end"
# Updates the Settings sheet's file-path / sheet-name constants to point at
# the new (Chile / Spain) survey shipment, and restores the scroll/selection
# state left by the author on Settings and Constants.

$wb = $excel.ActiveWorkbook

$settings  = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# --- Settings sheet value updates -----------------------------------------
$settings.Range("B2").Value  = "/Planeacion/0.Envios TS 2021/32 Envio Semana 32"
$settings.Range("B3").Value  = "REPORTE_ESTUDIOC_CHILE_SEM_49"
$settings.Range("B4").Value  = "TRADICIONAL"
$settings.Range("B5").Value  = "LUXURY"

$settings.Range("B7").Value  = "/Planeacion/0.Envios TS 2021/33 Envio Semana 31/ESPAÑA/LUXURY"
$settings.Range("B8").Value  = "BD_ENCUESTA_LARGA_GULF_VIP_Sema"

$settings.Range("B10").Value = "/Planeacion/0.Envios TS 2021/33 Envio Semana 31/ESPAÑA/LUXURY"
$settings.Range("B11").Value = "BD_ENCUESTA_LARGA_GULF_VIP_Sema"

$settings.Range("B13").Value = "/Planeacion/0.Envios TS 2021/33 Envio Semana 31/ESPAÑA/TRADICIONAL"
$settings.Range("B14").Value = "BD_ENCUESTA_LARGA_GULF_Tradicio"

$settings.Range("B16").Value = "/Planeacion/0.Envios TS 2021/32 Envio Semana 32/CHILE/Base de Datos"

$settings.Range("B19").Value = "/Planeacion/0.Envios TS 2021/32 Envio Semana 32/CHILE/Base de Datos/Exportadas"
$settings.Range("B22").Value = "/Planeacion/0.Envios TS 2021/32 Envio Semana 32/CHILE/Base de Datos/Exportadas"
$settings.Range("B25").Value = "/Planeacion/0.Envios TS 2021/32 Envio Semana 32/CHILE/Base de Datos/Exportadas"

$settings.Range("B28").Value = "/Planning/0.Shipping TS 2021/1 Directory"
$settings.Range("B29").Value = "CHILE"

$settings.Range("B36").Value = "/Planeacion/0.Envios TS 2021/1 Directorio"
$settings.Range("B37").Value = "España"

# B36 picks up the same (non-justified) look as B37 once its content becomes
# a plain path again instead of the old hyperlink-styled directory value.
$settings.Range("B37").Copy()
$settings.Range("B36").PasteSpecial(-4122)
$settings.Application.CutCopyMode = $false

# --- restore view state (scroll position + selection) ----------------------
$settings.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$settings.Range("B6").Select()

$constants.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$constants.Range("B8").Select()

$settings.Activate()
